# amend results in thesis and re-examine
#
# Re-examined survey data in columns G:O (rows 1-5): several raw scores
# were corrected, and the three "pm"/"pf"/"cm" condition labels (row 6)
# plus Excel's built-in Good/Bad/Neutral cell styles were added to the
# G:O block to flag the three measurement groups. All downstream
# AVERAGE/AVERAGEIF formulas recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Apply the built-in cell styles to the three column groups -----
# Call order matters: it fixes the order in which styles.xml allocates
# font/fill/cellXfs slots (1st call -> slot 1, 2nd -> slot 2, ...).
# M:O -> "Bad" (red), G:I -> "Neutral" (yellow), J:L -> "Good" (green)
$ws.Range("M1:O5").Style = "Bad"
$ws.Range("G1:I5").Style = "Neutral"
$ws.Range("J1:L5").Style = "Good"

# --- 2. Corrected raw values, row by row, columns G:O -----------------
# Row 1
$ws.Cells.Item(1, 7).Value  = 2   # G1
$ws.Cells.Item(1, 8).Value  = 3   # H1
$ws.Cells.Item(1, 9).Value  = 1   # I1
$ws.Cells.Item(1, 10).Value = 1   # J1
$ws.Cells.Item(1, 11).Value = 1   # K1
$ws.Cells.Item(1, 12).Value = 1   # L1
$ws.Cells.Item(1, 13).Value = 3   # M1
$ws.Cells.Item(1, 14).Value = 3   # N1
$ws.Cells.Item(1, 15).Value = 3   # O1

# Row 2
$ws.Cells.Item(2, 7).Value  = 1   # G2
$ws.Cells.Item(2, 8).Value  = 1   # H2
$ws.Cells.Item(2, 9).Value  = 1   # I2
$ws.Cells.Item(2, 10).Value = 1   # J2
$ws.Cells.Item(2, 11).Value = 1   # K2
$ws.Cells.Item(2, 12).Value = 1   # L2
$ws.Cells.Item(2, 13).Value = 1   # M2
$ws.Cells.Item(2, 14).Value = 1   # N2
$ws.Cells.Item(2, 15).Value = 1   # O2

# Row 3
$ws.Cells.Item(3, 7).Value  = 2   # G3
$ws.Cells.Item(3, 8).Value  = 0   # H3
$ws.Cells.Item(3, 9).Value  = 2   # I3
$ws.Cells.Item(3, 10).Value = 2   # J3
$ws.Cells.Item(3, 11).Value = 2   # K3
$ws.Cells.Item(3, 12).Value = 2   # L3
$ws.Cells.Item(3, 13).Value = 2   # M3
$ws.Cells.Item(3, 14).Value = 2   # N3
$ws.Cells.Item(3, 15).Value = 2   # O3

# Row 4
$ws.Cells.Item(4, 7).Value  = 3   # G4
$ws.Cells.Item(4, 8).Value  = 0   # H4
$ws.Cells.Item(4, 9).Value  = 0   # I4
$ws.Cells.Item(4, 10).Value = 3   # J4
$ws.Cells.Item(4, 11).Value = 3   # K4
$ws.Cells.Item(4, 12).Value = 3   # L4
$ws.Cells.Item(4, 13).Value = 3   # M4
$ws.Cells.Item(4, 14).Value = 3   # N4
$ws.Cells.Item(4, 15).Value = 3   # O4

# Row 5
$ws.Cells.Item(5, 7).Value  = 1   # G5
$ws.Cells.Item(5, 8).Value  = 0   # H5
$ws.Cells.Item(5, 9).Value  = 0   # I5
$ws.Cells.Item(5, 10).Value = 1   # J5
$ws.Cells.Item(5, 11).Value = 1   # K5
$ws.Cells.Item(5, 12).Value = 3   # L5
$ws.Cells.Item(5, 13).Value = 3   # M5
$ws.Cells.Item(5, 14).Value = 2   # N5
$ws.Cells.Item(5, 15).Value = 3   # O5

# --- 3. New label row 6 (shared-string labels for the three groups) ---
# Written in this order so sharedStrings.xml ends up as [cm, pm, pf].
$ws.Range("N6").Value = "cm"
$ws.Range("H6").Value = "pm"
$ws.Range("K6").Value = "pf"

# --- 4. Selection moved to D10 -----------------------------------------
$ws.Range("D10").Select()
